$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("device.json")

# 1. Insert a new row above row 5 ("macAddress" and everything below shifts down by one)
[void]$ws.Rows.Item(5).Insert()

# 2. Row 4 was "partNumber" with a Notes cell -> rename field to "board" and clear the note
$ws.Range("B4").Value = "board"
$ws.Range("D4").Value = ""

# 3. New row 5: "serialNumber" with a hyperlink in C5 to the IAO_0000131 term
$ws.Range("B5").Value = "serialNumber"
[void]$ws.Hyperlinks.Add($ws.Range("C5"), "http://purl.obolibrary.org/obo/IAO_0000131", "", "", "http://purl.obolibrary.org/obo/IAO_0000131")
$ws.Range("C5").Value = " http://purl.obolibrary.org/obo/IAO_0000131"

# 4. Row 6 (previously row 5, macAddress) now needs the semantic notation URL in C6
$ws.Range("C6").Value = "http://ns.cerise-project.nl/energy/def/cim-smartgrid/#ElectronicAddress.mac"

# 5. Resize Table1 so it covers the newly inserted row
$tbl = $ws.ListObjects.Item(1)
[void]$tbl.Resize($ws.Range("B2:D28"))

# 6. Widen column C to fit the new, longer content
$ws.Columns.Item(3).ColumnWidth = 71.27

# 7. Update the selection to match the target view
[void]$ws.Range("E10").Select()
